$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextDate($cellAddr, $text) {
    # Some "DD-MM-YYYY" strings (day <= 12) are ambiguous with MM-DD-YYYY and
    # Excel's smart-entry would silently reinterpret them as date serials.
    # Force text entry, then restore the "Normal" style so no stray
    # number-format/style ends up applied to the cell.
    $cell = $ws.Range($cellAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 3: date format + D/G flip to 1 (H stays 1)
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4: date format (ambiguous) + D/E to 1, H to 0
Set-TextDate "A4" "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5: date format only (ambiguous)
Set-TextDate "A5" "04-08-2022"

# Row 6: date format only (ambiguous)
Set-TextDate "A6" "08-08-2022"

# Row 7: date format only (ambiguous)
Set-TextDate "A7" "11-08-2022"

# Row 8: date format only
$ws.Range("A8").Value = "15-08-2022"

# Row 9: date format only
$ws.Range("A9").Value = "18-08-2022"

# Row 10: date format only
$ws.Range("A10").Value = "22-08-2022"

# Row 11: date format + D/E to 1, H to 0
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 0

# Row 12: date format + D/E to 1, H to 0
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

# Row 13: date format only (ambiguous)
Set-TextDate "A13" "01-09-2022"

# Row 14: date format (ambiguous) + D/E to 1, H to 0
Set-TextDate "A14" "05-09-2022"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("H14").Value = 0

# Row 15: date format only (ambiguous)
Set-TextDate "A15" "08-09-2022"

# Row 16: date format only (ambiguous)
Set-TextDate "A16" "12-09-2022"

# Row 17: date format only
$ws.Range("A17").Value = "15-09-2022"

# Row 18: date format only
$ws.Range("A18").Value = "19-09-2022"

# Row 19: date format only
$ws.Range("A19").Value = "22-09-2022"

# Row 20: date format + D/E to 1, H to 0
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("H20").Value = 0

# Row 21: date format only
$ws.Range("A21").Value = "29-09-2022"
